$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("K10").Value = 507.6
$ws1.Range("M14").Value = 4850.15
$ws1.Range("H16").Value = 970.08
$ws1.Range("M16").Value = 2186.02
$ws1.Range("M17").Value = 16276.37

$ws1.Range("H23").Value = "2 de 21"
$ws1.Range("I23").Value = "0 de 21"
$ws1.Range("J23").Value = "0 de 21"
$ws1.Range("K23").Value = "1 de 21"
$ws1.Range("M23").Value = "5 de 21"

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F10").Value = 507.6
$ws2.Range("F14").Value = 4850.15
$ws2.Range("F16").Value = 3156.1
$ws2.Range("F17").Value = 16276.37
$ws2.Range("F23").Value = 43125.96

$ws2.Columns.Item(6).ColumnWidth = 13.16666666666667

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws3.Range("D6").Value = 2828.71
$ws3.Range("E6").Value = 78.87368146025983
$ws3.Range("F6").Value = 0.972873117302458

$ws3.Range("D10").Value = 507.6
$ws3.Range("E10").Value = -119.492016465608
$ws3.Range("F10").Value = 1.307883428156843

$ws3.Range("D12").Value = 35759.6
$ws3.Range("E12").Value = 1980.139999999999
$ws3.Range("F12").Value = 0.9475317000064124

$ws3.Range("D14").Value = 43125.96
$ws3.Range("E14").Value = 12298.78147880389
$ws3.Range("F14").Value = 0.7780994344645284

$ws3.Columns.Item(4).ColumnWidth = 12.16666666666667
$ws3.Columns.Item(5).ColumnWidth = 22.16666666666667
